# Notifications.xlsx - "Added 2 new watchlist test cases"
#
# The "TCID=Notifications0008" row (row 9) on the "Test Cases" sheet had its
# Jira-id / Description cell contents edited: one item was dropped from each
# of the two pipe-delimited ("||") lists (the "...published a comment on
# post..." test case and its matching "OPQA-1397" id were removed), and the
# row height shrank to match the now-shorter wrapped text. The sheet's
# active-cell selection also moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Column B = "Jira id", Column C = "Description" (row 9)
$ws.Range("B9").Value = "OPQA-877||OPQA-1013||OPQA-215||OPQA-1395"
$ws.Range("C9").Value = "Verify that user receives a notification when someone he is following  publishes a post||Verify that user is receiving notification when someone liked his post(aggregated notification)||Verify that user able to recevies a notification when other user commented on his post||Verify that all users receive notification when other user published a post and validate notification."

# Row 9 now wraps to a shorter height after the text trim.
$ws.Rows.Item(9).RowHeight = 45

# Selection moved from C8 to C13.
$ws.Range("C13").Select()
